$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("title + description")

# Activate the sheet (it is already tabSelected) and set the new view
$ws.Activate()

# New row of data for the web-design "under construction" placeholder page
$ws.Range("B20").Value = "webdesign.html"
$ws.Range("C20").Value = "Створення сайтів. Ми працюємо над цією сторінкою -   IT майстерня “Все працює”"
$ws.Range("D20").Value = "<meta name=""description"" content=""&#128736; Сторінка створюється.&#10057;   "">"
$ws.Rows.Item(20).RowHeight = 30

# Update the selection / view position to match the new edit location
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("D20").Select() | Out-Null
